$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:E hold numeric-looking text (e.g. "0.8447", "208"); force text format
# so Excel keeps them as strings rather than auto-converting to numbers
# (matches the t="inlineStr"/text storage used in the source file).
$ws.Range("B2:E11").NumberFormat = "@"

$data = @(
    @("A",   "0.8447", "0.8894", "0.8665", "208"),
    @("AA",  "0.8824", "0.5769", "0.6977", "52"),
    @("AAA", "0.9048", "0.7917", "0.8444", "24"),
    @("B",   "0.9272", "0.9091", "0.9180", "154"),
    @("BB",  "0.9357", "0.9225", "0.9291", "284"),
    @("BBB", "0.9129", "0.9532", "0.9326", "363"),
    @("C",   "1.0000", "1.0000", "1.0000", "4"),
    @("CC",  "1.0000", "0.5000", "0.6667", "2"),
    @("CCC", "0.7857", "0.8462", "0.8148", "26"),
    @("D",   "1.0000", "1.0000", "1.0000", "1")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
